$d = $word.ActiveDocument

# 1. Update activation date: 2012 -> 2023
$d.Content.Find.Execute("Ativação: 01/01/2012", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Ativação: 01/01/2023", 2)

# Helper: find the paragraph whose text equals $marker (ignoring trailing
# paragraph/line-break marks) and insert a brand new paragraph right after
# it, containing $newText in italics - mirroring what Word does when you
# place the cursor at the end of a paragraph, press Enter and type new
# (italic) text.
function Insert-ItalicParagraphAfter($marker, $newText) {
    $doc = $word.ActiveDocument
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $marker) {
            $p.Range.InsertParagraphAfter()
            $newPara = $doc.Paragraphs($p.Index + 1)
            $r = $newPara.Range
            $start = $r.Start
            $r.InsertAfter($newText)
            $target = $doc.Range($start, $start + $newText.Length)
            $target.Font.Italic = $true
            break
        }
    }
}

# 2. New italic English sentence after the "Objetivos" paragraph
Insert-ItalicParagraphAfter `
    "Fornecer ao aluno seminários sobre temas atuais de Física, Tecnologia e Engenharia." `
    "Provide student seminars on current topics in Physics, Technology and Engineering."

# 3. New italic English sentence after the "Programa resumido" paragraph
Insert-ItalicParagraphAfter `
    "Seminários abrangendo os cenários atuais e futuros da indústria de alta tecnologia e do campo de atuação do engenheiro físico." `
    "Seminars covering the current and future scenarios of the high technology industry and the field of activity of the physical engineer."

# 4. New italic English sentence after the "Programa" paragraph
Insert-ItalicParagraphAfter `
    "Seminários seguido de debates com profissionais e estudantes de graduação e pós-graduação sobre temas relevantes e atuais das áreas de Física, Tecnologia e Engenharia, abrangendo desde as pesquisas básicas até o segmento industrial e de serviços." `
    "Seminars followed by debates with professionals and undergraduate and graduate students on relevant and current topics in the areas of Physics, Technology and Engineering, ranging from basic research to the industrial and services segment."

# 5. Drop the " I" after "Projeto Integrado" in the prerequisite line
$d.Content.Find.Execute("LOM3238 -  Projeto Integrado I  (Requisito)", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "LOM3238 -  Projeto Integrado  (Requisito)", 2)
